# Updated symbol list on Thu Feb  2 04:21:38 UTC 2023 with GitHub Actions
#
# Refresh cryptocurrency price/volume snapshot on Sheet1. Columns:
#   B = Coin name, C = coinranking.com link, D = Price, E = Volume(1h)
# All D/E values are stored as text (matching the source feed's inline
# strings), so numeric-looking values are entered with a leading
# apostrophe to force text and the cell style is reset back to "Normal"
# afterwards so no stray number-format is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'326.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.65%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.254"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.19%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08092"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.25%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.93%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.912"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.00%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-1.39%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9367"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.02%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1325"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'17.98%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1954"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.24%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09173"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.40%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03436"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.05%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09543"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.59%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001398"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.50%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = "'0.006003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'2.22%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = "'3.360"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-6.47%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = "'0.3524"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'3.38%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").Value = "'7.249"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'22.53%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = "'0.1314"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.66%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").Value = "'0.2313"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-10.66%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").Value = "'0.04446"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.77%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.96%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004352"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-6.13%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001290"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-5.21%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'0.02%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02469"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'10.36%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05240"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.48%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007694"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.02%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1430"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.51%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008624"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-4.39%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002161"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.33%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008167"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-5.54%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.18%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E48").Value = "'-13.11%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'148.12%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
